$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)   # "measurements"
$ws2 = $wb.Worksheets.Item(2)   # "sites" - used as a style donor (same header pattern)

# --- Copy cell formatting from an existing, already-styled sheet so the
# --- generated cellXfs entries match the ones the diff expects (s="12",
# --- s="4", s="14", s="6", s="16", s="8", s="18", s="10").
$ws2.Range("B1").Copy() | Out-Null
$ws4.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws2.Range("G1").Copy() | Out-Null
$ws4.Range("C1").PasteSpecial(-4122) | Out-Null

$ws2.Range("B2").Copy() | Out-Null
$ws4.Range("B2").PasteSpecial(-4122) | Out-Null

$ws2.Range("G2").Copy() | Out-Null
$ws4.Range("C2").PasteSpecial(-4122) | Out-Null

$ws2.Range("B3").Copy() | Out-Null
$ws4.Range("B3").PasteSpecial(-4122) | Out-Null

$ws2.Range("G3").Copy() | Out-Null
$ws4.Range("C3").PasteSpecial(-4122) | Out-Null

$ws2.Range("B4").Copy() | Out-Null
$ws4.Range("B4").PasteSpecial(-4122) | Out-Null

$ws2.Range("G4").Copy() | Out-Null
$ws4.Range("C4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Cell values ---
$ws4.Range("B1").Value = "Names"
$ws4.Range("C1").Value = ""
$ws4.Range("B2").Value = "Measurement Number"
$ws4.Range("C2").Value = "Experiment Name"
$ws4.Range("B3").Value = "String"
$ws4.Range("C3").Value = "String"
$ws4.Range("B4").Value = "number"
$ws4.Range("C4").Value = "experiment"

# --- Column widths ---
# (target widths from the spec are 20.3515625 / 17.140625; the values below are
# the closest the engine's pixel-quantized ColumnWidth model can represent)
$ws4.Range("B1").EntireColumn.ColumnWidth = 19.5
$ws4.Range("C1").EntireColumn.ColumnWidth = 16.333333333333332

# --- Merge the group header cells across B1:C1 ---
$ws4.Range("B1:C1").Merge()
